# "updated on 26 Mar 2018"
# AutomationControlSheet.xlsx settings update:
#   - OS / BROWSER defaults changed
#   - CLOUD PROVIDER contact e-mail changed and turned into a mailto: link
#   - the URL and REST-API host cells turned into (or restyled as) hyperlinks
#   - a new conditional-formatting rule greys out the OS cell when running Local
#   - two new Project/URL lookup rows added to the "values" sheet, each linked

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AppControl")
$ws2 = $wb.Worksheets.Item("values")

# ---------------------------------------------------------------------------
# AppControl sheet - plain value edits
# ---------------------------------------------------------------------------
$ws1.Range("B12").Value = "WINDOWS_10"   # OS
$ws1.Range("B13").Value = "FIREFOX"      # BROWSER

# ---------------------------------------------------------------------------
# AppControl sheet - hyperlink / style changes
# ---------------------------------------------------------------------------

# CLOUD PROVIDER: new e-mail address, becomes a mailto hyperlink (no border)
$ws1.Hyperlinks.Add($ws1.Range("B15"), "mailto:raghav.pal@testingxperts.com", "", "", "raghav.pal@testingxperts.com")

# URL (B6) already behaved like a hyperlink - it now loses its cell border
$ws1.Range("B6").Borders.LineStyle = -4142   # xlLineStyleNone

# REST API "URL OR HOSTNAME" becomes a real hyperlink and keeps a border,
# matching the other bordered hyperlink cells (B30:B32)
$ws1.Hyperlinks.Add($ws1.Range("B40"), "https://townsqd.com/api", "", "", "https://townsqd.com/api")
$ws1.Range("B40").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# AppControl sheet - new conditional formatting: grey out OS (B12) on Local
# ---------------------------------------------------------------------------
$cf = $ws1.Range("B12").FormatConditions.Add(2, 0, '=$B$8="Local"')
$cf.Interior.Color = 8421504   # RGB(128,128,128) == theme0 tint -0.499984740745262

# ---------------------------------------------------------------------------
# AppControl sheet - selection / viewport
# ---------------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("B44").Select()

# ---------------------------------------------------------------------------
# values sheet - two new Project -> URL lookup rows, each hyperlinked
# ---------------------------------------------------------------------------
$ws2.Range("B12").Value = "UPC"
$ws2.Range("C12").Value = "https://uat-upcenhancements.hostedinsurance.com/AgentPortal/login"
$ws2.Hyperlinks.Add($ws2.Range("C12"), "https://uat-upcenhancements.hostedinsurance.com/AgentPortal/login", "", "", "https://uat-upcenhancements.hostedinsurance.com/AgentPortal/login")

$ws2.Range("B13").Value = "VitaminShoppe"
$ws2.Range("C13").Value = "https://redesign.perf.vitaminshoppe.com/"
$ws2.Hyperlinks.Add($ws2.Range("C13"), "https://redesign.perf.vitaminshoppe.com/", "", "", "https://redesign.perf.vitaminshoppe.com/")
$ws2.Range("C13").Borders.LineStyle = 1
